{"js": "// \"Ch\u1ec9nh l\u1ea1i m\u1eabu 26\" \u2014 remove the leftover \"vnpt.SiteAddress\" merge-field\n// placeholder text that trails the \"\u0110\u1ecba ch\u1ec9: \" label for B\u00ean A, leaving the\n// label run (\"\u0110\u1ecba ch\u1ec9: \") untouched in its bulleted paragraph.\n\nconst body = context.document.body;\nconst results = body.search(\"vnpt.SiteAddress\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# \"Ch\u1ec9nh l\u1ea1i m\u1eabu 26\" \u2014 remove the leftover \"vnpt.SiteAddress\" merge-field\n# placeholder text that trails the \"\u0110\u1ecba ch\u1ec9: \" label for B\u00ean A, leaving the\n# label run (\"\u0110\u1ecba ch\u1ec9: \") untouched in its bulleted paragraph.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"vnpt.SiteAddress\"\n$range.Find.MatchCase = $true\n$range.Find.Wrap = 1\n\nwhile ($range.Find.Execute()) {\n    $range.Delete()\n}\n"}
